# Add new equipment-log entries (row 17, row 18) to the "main" sheet, push the
# trailing "**" marker down to row 19, and add a new "average cycles between
# replacements" summary column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# --- New log entry: 07/02/2018, operator 3012, cracked part, 1001 strokes ---
$ws.Range("A17").Value = "07/02/2018"
$ws.Range("B17").Value = "3012"
$ws.Range("C17").Value = "Тріснута запчастина"
$ws.Range("D17").Value = 1001

# --- New log entry: 08/02/2018, operator 3333, contact asymmetry, 5698 strokes ---
$ws.Range("A18").Value = "08/02/2018"
$ws.Range("B18").Value = "3333"
$ws.Range("C18").Value = "Асиметрія контакту"
$ws.Range("D18").Value = 5698

# --- Move the trailing "**" separator marker down from row 17 to row 19 ---
$ws.Range("A19").Value = "**"

# --- New summary column H: average number of cycles between replacements ---
$ws.Range("H1").Value = "середне значення циклів між замінами"
$ws.Range("H2").Value = 150000
